$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new data row right below the header row ---------------------
# New most-recent record for the first CEDULA group (0105781496); pushes all
# existing data rows down by one (old row 2 -> row 3, etc.).
$ws.Rows.Item(2).Insert()

$ws.Cells.Item(2, 1).Value = 8

# Column B holds a CEDULA stored as text (leading zeros matter). Assigning a
# numeric-looking string via .Value auto-converts it to a number, so mark the
# cell as text first (quote-prefix), then restore the default "Normal" style
# so no extra number-format / quote-prefix metadata is left behind.
$ws.Cells.Item(2, 2).Value = "'0105781496"
$ws.Cells.Item(2, 2).Style = "Normal"

# Column C is a date/time serial; give it the same custom date format used by
# the rest of the column so it reuses the existing style instead of staying
# General.
$ws.Cells.Item(2, 3).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(2, 3).Value = 45257.83576388889

$ws.Cells.Item(2, 4).Value = "SwR2Wh"

# --- Insert a new data row before the final two rows -----------------------
# New most-recent record for the second CEDULA group (0123456789); pushes the
# remaining old rows down by one (old row 7 -> row 8, old row 8 -> row 9).
$ws.Rows.Item(8).Insert()

$ws.Cells.Item(8, 1).Value = 9

$ws.Cells.Item(8, 2).Value = "'0123456789"
$ws.Cells.Item(8, 2).Style = "Normal"

$ws.Cells.Item(8, 3).NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Cells.Item(8, 3).Value = 45257.83667824074

$ws.Cells.Item(8, 4).Value = "3D9AlD"
